$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Chris"
$ws.Range("C11").Value = 432594785
$ws.Range("D11").Value = "house"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "9 Hetherington"
$ws.Range("I11").Value = "Herston"
$ws.Range("J11").Value = 4006
$ws.Range("K11").Value = 380

$ws.Range("A12").Select()
